$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two data rows were removed entirely from the table (missing-data
# re-sampling), which shifts every following row up.
#   Row 26 = "RM 232"
#   Row 28 = "SC 92" (becomes row 27 after the first delete)
$ws.Rows(26).Delete()
$ws.Rows(26).Delete()

# Remaining per-cell value changes (some cells that had data are now
# blanked out / marked missing, some previously-missing cells now have
# an imputed/observed value).
$ws.Range("C2").Value = 14.9
$ws.Range("C3").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("F6").Value = 16.43
$ws.Range("D8").Value = -13.9
$ws.Range("D10").Value = -14.7
$ws.Range("C11").Value = 11.4
$ws.Range("F11").Value = 17.65
$ws.Range("D12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("F13").Value = 17.1
$ws.Range("D15").Value = -15.2
$ws.Range("F17").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("F18").Value = 18.35
$ws.Range("D19").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("C21").Value = 12.7
$ws.Range("F24").Value = ""
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = -15.5
$ws.Range("F25").Value = 16.6
$ws.Range("D27").Value = -14.6
$ws.Range("B29").Value = ""
$ws.Range("D29").Value = ""
$ws.Range("F31").Value = ""
$ws.Range("F32").Value = ""
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = ""
